$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Set D4 mark value (Session 3, D&C) for the student
$ws.Range("D4").Value = 7

# Add grading comment for Session 3 (D&C) in D5, mirroring the existing
# comment style used for other sessions
$comment = "The D&C version is not working as expected. It works but the complexity is not O(nlogn). Please, check the video of the last seminar (the implementation should be very similar to Mergesort)." + [char]10 + "Nevative number of inversions are because you should use long instead of int (there are very large results)."
$ws.Range("D5").Value = $comment

# Writing the long text can trigger an automatic row-height adjustment on
# row 5 (wrap text style); keep it at its original height.
$ws.Rows.Item(5).RowHeight = 12.75

# Increase the height of row 12 to fit the longer comment text
$ws.Rows.Item(12).RowHeight = 83.25

# Update the selected range to reflect where the user was working
$ws.Range("D5:D12").Select()
